$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.527.93"

$ws.Range("D3").Value = "1.913.36"

$ws.Range("E3").Value = "  +4.43%  "

$ws.Range("E4").Value = "  +0.22%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "315.05"
$c.Style = "Normal"

$ws.Range("E5").Value = "  +1.53%  "

$ws.Range("E6").Value = "  +0.06%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5149"
$c.Style = "Normal"

$ws.Range("E7").Value = "  +3.79%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3978"
$c.Style = "Normal"

$ws.Range("E8").Value = "  +1.29%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.09825"
$c.Style = "Normal"

$ws.Range("E9").Value = "  -4.00%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.152"
$c.Style = "Normal"

$ws.Range("E10").Value = "  +3.50%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "42.26"
$c.Style = "Normal"

$ws.Range("E11").Value = "  +2.68%  "

$ws.Range("E12").Value = "  +1.51%  "

$ws.Range("E13").Value = "  +2.45%  "

$ws.Range("D14").Value = "1.908.67"

$ws.Range("E14").Value = "  +4.61%  "

$ws.Range("E15").Value = "  +1.92%  "

$ws.Range("E16").Value = "  +0.16%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "94.71"
$c.Style = "Normal"

$ws.Range("E17").Value = "  +1.79%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001137"
$c.Style = "Normal"

$ws.Range("E18").Value = "  -1.32%  "

$ws.Range("E19").Value = "  +0.10%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "18.28"
$c.Style = "Normal"

$ws.Range("E20").Value = "  +5.89%  "

$ws.Range("E21").Value = "  +0.09%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.329"
$c.Style = "Normal"

$ws.Range("E22").Value = "  +5.27%  "

$ws.Range("D23").Value = "28.585.66"

$ws.Range("E23").Value = "  +1.16%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.51"
$c.Style = "Normal"

$ws.Range("E24").Value = "  +1.86%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.321"
$c.Style = "Normal"

$ws.Range("E25").Value = "  +3.07%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.685"
$c.Style = "Normal"

$ws.Range("E26").Value = "  +10.23%  "

$ws.Range("D27").Value = "2.129.15"

$ws.Range("E27").Value = "  +4.47%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "21.30"
$c.Style = "Normal"

$ws.Range("E28").Value = "  +2.41%  "

$ws.Range("E29").Value = "  -0.19%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "129.05"
$c.Style = "Normal"

$ws.Range("E30").Value = "  +1.89%  "

$ws.Range("E31").Value = "  +7.31%  "

$ws.Range("E32").Value = "  +2.56%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.755"
$c.Style = "Normal"

$ws.Range("E33").Value = "  +2.77%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.630"
$c.Style = "Normal"

$ws.Range("E34").Value = "  +0.84%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "9.872"
$c.Style = "Normal"

$ws.Range("E35").Value = "  +9.27%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.06796"
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02437"
$c.Style = "Normal"

$ws.Range("E37").Value = "  +3.11%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.272"
$c.Style = "Normal"

$ws.Range("E38").Value = "  +6.82%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.2216"
$c.Style = "Normal"

$ws.Range("E39").Value = "  +2.94%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "11.82"
$c.Style = "Normal"

$ws.Range("E40").Value = "  +3.23%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.6488"
$c.Style = "Normal"

$ws.Range("E41").Value = "  +4.18%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.096"
$c.Style = "Normal"

$ws.Range("E42").Value = "  +2.11%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.189"
$c.Style = "Normal"

$ws.Range("E43").Value = "  +1.38%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("E45").Value = "  +2.71%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.6113"
$c.Style = "Normal"

$ws.Range("E46").Value = "  +2.86%  "

$ws.Range("E47").Value = "  +2.39%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.292"
$c.Style = "Normal"

$ws.Range("E48").Value = "  +1.65%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.047"
$c.Style = "Normal"

$ws.Range("E49").Value = "  +4.81%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "125.00"
$c.Style = "Normal"

$ws.Range("E50").Value = "  +0.68%  "

$ws.Range("E51").Value = "  +1.80%  "
